# Update the "Translation" sheet rows 5-14 (TEXT ID / ALIGNMENT / GB columns)
# to reflect the re-shuffled/renamed translation entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @{ Row = 5;  B = "SingleUseId3";  D = "Center"; F = "btn2" },
    @{ Row = 6;  B = "SingleUseId11"; D = "Left";   F = "page 2" },
    @{ Row = 7;  B = "SingleUseId12"; D = "Left";   F = "TouchGFX TEST" },
    @{ Row = 8;  B = "SingleUseId2";  D = "Center"; F = "btn1" },
    @{ Row = 9;  B = "SingleUseId10"; D = "Left";   F = "page 1" },
    @{ Row = 10; B = "SingleUseId9";  D = "Left";   F = "TouchGFX TEST" },
    @{ Row = 11; B = "SingleUseId5";  D = "Center"; F = "1" },
    @{ Row = 12; B = "SingleUseId6";  D = "Center"; F = "4" },
    @{ Row = 13; B = "SingleUseId7";  D = "Center"; F = "3" },
    @{ Row = 14; B = "SingleUseId8";  D = "Center"; F = "2" }
)

# A scratch cell used to push purely-numeric-looking text ("1","2","3","4")
# into a cell as TEXT rather than as a Number. Setting .Value directly with
# a digit string makes Excel coerce it into a numeric cell (no t="s"); going
# through a TEXT() formula + copy/paste-values keeps the literal a string,
# matching the original workbook's shared-string-backed text cells and
# without touching any cell's NumberFormat/style.
$scratch = $ws.Range("Z1")

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("D" + $r.Row).Value = $r.D

    $fCell = $ws.Range("F" + $r.Row)
    if ($r.F -match '^[0-9]+$') {
        $scratch.Formula = '=TEXT(' + $r.F + ',"0")'
        $scratch.Copy()
        $fCell.PasteSpecial(-4163)
    } else {
        $fCell.Value = $r.F
    }
}

$scratch.ClearContents()
